# Model export fixes + slider battery and slider windfarm addition to generic
# excel templates. This adds a new "SLIDER_WF" template row (row 2) to the
# windfarms sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for the slider windfarm template entry.
$ws.Range("A2").Value = "SLIDER_WF"
$ws.Range("B2").Value = "SLIDER_WF"
$ws.Range("C2").Value = "SLIDER_WF"
$ws.Range("K2").Value = $true
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 5000
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 5000
$ws.Range("P2").Value = 52
$ws.Range("Q2").Value = 5

# Match the saved selection/active cell from the authoring session.
$ws.Range("N12").Select() | Out-Null
